$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1418.7717
$ws.Range("I15").Value = 1418.7717
$ws.Range("K15").Value = 4256.3151
$ws.Range("M15").Value = -4087.3151
$ws.Range("H17").Value = 1450.338
$ws.Range("J17").Value = 1450.338
$ws.Range("L17").Value = 4351.014
$ws.Range("N17").Value = -4687.014
$ws.Range("H29").Value = 1060
$ws.Range("I29").Value = 850
$ws.Range("K29").Value = 2550
$ws.Range("M29").Value = -2269
$ws.Range("H133").Value = 77660
$ws.Range("J133").Value = 77660
$ws.Range("L133").Value = 77660
$ws.Range("N133").Value = -87780
$ws.Range("H134").Value = 115992
$ws.Range("J134").Value = 115992
$ws.Range("L134").Value = 115992
$ws.Range("N134").Value = -126132
$ws.Range("H138").Value = 3338344.8
$ws.Range("I138").Value = 8698603
$ws.Range("J138").Value = 6292.1353
$ws.Range("K138").Value = 26095809
$ws.Range("L138").Value = 18876.4059
$ws.Range("M138").Value = -26090669
$ws.Range("N138").Value = -29156.4059
$ws.Range("H139").Value = 80780
$ws.Range("J139").Value = 80780
$ws.Range("L139").Value = 80780
$ws.Range("N139").Value = -91060
$ws.Range("H140").Value = 79433.84
$ws.Range("J140").Value = 78257.5
$ws.Range("L140").Value = 78257.5
$ws.Range("N140").Value = -88617.5
$ws.Range("H141").Value = 22127
$ws.Range("I141").Value = 17918
$ws.Range("J141").Value = 26336
$ws.Range("K141").Value = 53754
$ws.Range("L141").Value = 79008
$ws.Range("M141").Value = -48574
$ws.Range("N141").Value = -89368

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14189.339
$ws.Range("I32").Value = 12538.362
$ws.Range("J32").Value = 38128.5
$ws.Range("K32").Value = 12538.362
$ws.Range("L32").Value = 38128.5
$ws.Range("M32").Value = -12251.362
$ws.Range("N32").Value = -38702.5
$ws.Range("H63").Value = 10931.125
$ws.Range("I63").Value = 13326.667
$ws.Range("J63").Value = 9493.799999999999
$ws.Range("K63").Value = 13326.667
$ws.Range("L63").Value = 9493.799999999999
$ws.Range("M63").Value = -12640.667
$ws.Range("N63").Value = -10865.8
$ws.Range("H66").Value = 10931.125
$ws.Range("I66").Value = 13326.667
$ws.Range("J66").Value = 9493.799999999999
$ws.Range("K66").Value = 66633.33499999999
$ws.Range("L66").Value = 47469
$ws.Range("M66").Value = -63201.33499999999
$ws.Range("N66").Value = -54333
$ws.Range("H132").Value = 771401.7
$ws.Range("I132").Value = 953945.0600000001
$ws.Range("J132").Value = 4719.6
$ws.Range("K132").Value = 2861835.18
$ws.Range("L132").Value = 14158.8
$ws.Range("M132").Value = -2859305.18
$ws.Range("N132").Value = -19218.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3546.2
$ws.Range("I105").Value = 3085.2144
$ws.Range("K105").Value = 3085.2144
$ws.Range("M105").Value = -1338.2144
$ws.Range("H134").Value = 437481.7
$ws.Range("I134").Value = 608289.4
$ws.Range("J134").Value = 3892.923
$ws.Range("K134").Value = 1824868.2
$ws.Range("L134").Value = 11678.769
$ws.Range("M134").Value = -1822333.2
$ws.Range("N134").Value = -16748.769
$ws.Range("H140").Value = 49622.145
$ws.Range("J140").Value = 49622.145
$ws.Range("L140").Value = 49622.145
$ws.Range("N140").Value = -59982.145

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1950812.4
$ws.Range("I58").Value = 2470689
$ws.Range("K58").Value = 2470689
$ws.Range("M58").Value = -2470486
$ws.Range("H62").Value = 76622.86
$ws.Range("J62").Value = 8950
$ws.Range("L62").Value = 8950
$ws.Range("N62").Value = -10198
$ws.Range("H65").Value = 76622.86
$ws.Range("J65").Value = 8950
$ws.Range("L65").Value = 44750
$ws.Range("N65").Value = -50990
$ws.Range("H129").Value = 44634.223
$ws.Range("J129").Value = 44634.223
$ws.Range("L129").Value = 44634.223
$ws.Range("N129").Value = -54634.223
$ws.Range("H132").Value = 543098.5600000001
$ws.Range("I132").Value = 846919.5600000001
$ws.Range("J132").Value = 2972.3333
$ws.Range("K132").Value = 2540758.68
$ws.Range("L132").Value = 8916.999899999999
$ws.Range("M132").Value = -2538228.68
$ws.Range("N132").Value = -13976.9999
$ws.Range("H134").Value = 1850.9474
$ws.Range("I134").Value = 1633.5
$ws.Range("J134").Value = 2459.8
$ws.Range("K134").Value = 4900.5
$ws.Range("L134").Value = 7379.400000000001
$ws.Range("M134").Value = -2365.5
$ws.Range("N134").Value = -12449.4
$ws.Range("H136").Value = 1950812.4
$ws.Range("I136").Value = 2470689
$ws.Range("K136").Value = 7412067
$ws.Range("M136").Value = -7409517

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2927.3333
$ws.Range("I5").Value = 5121.6
$ws.Range("J5").Value = 1360
$ws.Range("K5").Value = 15364.8
$ws.Range("L5").Value = 4080
$ws.Range("M5").Value = -15252.8
$ws.Range("N5").Value = -4304
$ws.Range("H32").Value = 1199.8572
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 1199.8572
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 3599.5716
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -4165.571599999999
$ws.Range("H39").Value = 4397.1
$ws.Range("J39").Value = 4397.1
$ws.Range("L39").Value = 13191.3
$ws.Range("N39").Value = -13779.3
$ws.Range("H40").Value = 866.6667
$ws.Range("I40").Value = 56.857143
$ws.Range("J40").Value = 1575.25
$ws.Range("K40").Value = 227.428572
$ws.Range("L40").Value = 6301
$ws.Range("M40").Value = -158.428572
$ws.Range("N40").Value = -6439
$ws.Range("H46").Value = 3133.6316
$ws.Range("I46").Value = 256.66666
$ws.Range("J46").Value = 4461.4614
$ws.Range("K46").Value = 769.9999799999999
$ws.Range("L46").Value = 13384.3842
$ws.Range("M46").Value = -678.9999799999999
$ws.Range("N46").Value = -13566.3842
$ws.Range("H57").Value = 3512.5
$ws.Range("J57").Value = 3871.4285
$ws.Range("L57").Value = 11614.2855
$ws.Range("N57").Value = -12732.2855
$ws.Range("H58").Value = 2500
$ws.Range("I58").Value = 1000
$ws.Range("J58").Value = 2833.3333
$ws.Range("K58").Value = 3000
$ws.Range("L58").Value = 8499.999899999999
$ws.Range("M58").Value = -2872
$ws.Range("N58").Value = -8755.999899999999
$ws.Range("H135").Value = 2927.3333
$ws.Range("I135").Value = 5121.6
$ws.Range("J135").Value = 1360
$ws.Range("K135").Value = 46094.4
$ws.Range("L135").Value = 12240
$ws.Range("M135").Value = -43559.4
$ws.Range("N135").Value = -17310
$ws.Range("H138").Value = 1979.125
$ws.Range("J138").Value = 3501
$ws.Range("L138").Value = 10503
$ws.Range("N138").Value = -20783
$ws.Range("H139").Value = 1665.4348
$ws.Range("I139").Value = 1291.0555
$ws.Range("J139").Value = 3013.2
$ws.Range("K139").Value = 3873.1665
$ws.Range("L139").Value = 9039.599999999999
$ws.Range("M139").Value = 1266.8335
$ws.Range("N139").Value = -19319.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5710.967
$ws.Range("I70").Value = 4858.857
$ws.Range("J70").Value = 6456.5625
$ws.Range("K70").Value = 4858.857
$ws.Range("L70").Value = 6456.5625
$ws.Range("M70").Value = -4588.857
$ws.Range("N70").Value = -6996.5625
$ws.Range("H73").Value = 5710.967
$ws.Range("I73").Value = 4858.857
$ws.Range("J73").Value = 6456.5625
$ws.Range("K73").Value = 4858.857
$ws.Range("L73").Value = 6456.5625
$ws.Range("M73").Value = -3922.857
$ws.Range("N73").Value = -8328.5625
$ws.Range("H80").Value = 3389.4443
$ws.Range("I80").Value = 3000.7144
$ws.Range("K80").Value = 3000.7144
$ws.Range("M80").Value = -2002.7144
$ws.Range("H83").Value = 3389.4443
$ws.Range("I83").Value = 3000.7144
$ws.Range("K83").Value = 15003.572
$ws.Range("M83").Value = -10011.572
$ws.Range("H109").Value = 30282
$ws.Range("J109").Value = 30282
$ws.Range("L109").Value = 30282
$ws.Range("N109").Value = -32362
$ws.Range("H123").Value = 9968.134
$ws.Range("J123").Value = 9968.134
$ws.Range("L123").Value = 9968.134
$ws.Range("N123").Value = -14868.134
$ws.Range("H133").Value = 62320
$ws.Range("J133").Value = 62320
$ws.Range("L133").Value = 62320
$ws.Range("N133").Value = -72440
$ws.Range("H136").Value = 29628.572
$ws.Range("J136").Value = 29628.572
$ws.Range("L136").Value = 88885.716
$ws.Range("N136").Value = -93985.716
$ws.Range("H140").Value = 49687.5
$ws.Range("J140").Value = 49687.5
$ws.Range("L140").Value = 49687.5
$ws.Range("N140").Value = -60047.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 8600
$ws.Range("I100").Value = 11960
$ws.Range("K100").Value = 11960
$ws.Range("M100").Value = -11419
$ws.Range("H132").Value = 8703.412
$ws.Range("I132").Value = 8765.615
$ws.Range("J132").Value = 8501.25
$ws.Range("K132").Value = 26296.845
$ws.Range("L132").Value = 25503.75
$ws.Range("M132").Value = -23766.845
$ws.Range("N132").Value = -30563.75
$ws.Range("H137").Value = 87220
$ws.Range("J137").Value = 87220
$ws.Range("L137").Value = 87220
$ws.Range("N137").Value = -97420
$ws.Range("H139").Value = 48140
$ws.Range("J139").Value = 48140
$ws.Range("L139").Value = 48140
$ws.Range("N139").Value = -58420

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 23636.295
$ws.Range("J123").Value = 23636.295
$ws.Range("L123").Value = 23636.295
$ws.Range("N123").Value = -33436.295
$ws.Range("H125").Value = 70715
$ws.Range("J125").Value = 70715
$ws.Range("L125").Value = 70715
$ws.Range("N125").Value = -80555
$ws.Range("H132").Value = 2385.5417
$ws.Range("I132").Value = 1879.9412
$ws.Range("K132").Value = 5639.8236
$ws.Range("M132").Value = -3109.8236
$ws.Range("H139").Value = 60716.25
$ws.Range("J139").Value = 60716.25
$ws.Range("L139").Value = 60716.25
$ws.Range("N139").Value = -70996.25
